$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Profession sheet:
#  - "Profession" (B2) renamed to "Designation"
#  - "HomeTown" column removed (was H) -- remaining columns shift left
#  - "Migrated" (now at K after the shift) renamed to "MigratedFromHomeTown"
#  - Four new trailing columns added: FatherOccupation, MotherOccupation,
#    BrothersOccupation, SistersOccupation
# ---------------------------------------------------------------------------
$wsProfession = $wb.Worksheets.Item("Profession")
$wsProfession.Range("B2").Value = "Designation"
$wsProfession.Range("H2").EntireColumn.Delete()
$wsProfession.Range("K2").Value = "MigratedFromHomeTown"
$wsProfession.Range("N2").Value = "FatherOccupation"
$wsProfession.Range("O2").Value = "MotherOccupation"
$wsProfession.Range("P2").Value = "BrothersOccupation"
$wsProfession.Range("Q2").Value = "SistersOccupation"

# ---------------------------------------------------------------------------
# Property sheet:
#  - "SettledProperty" column (last column, O2) removed
# ---------------------------------------------------------------------------
$wsProperty = $wb.Worksheets.Item("Property")
$wsProperty.Range("O2").ClearContents()

# ---------------------------------------------------------------------------
# Education sheet:
#  - "HighestEducation" (B2) renamed to "HighestQualification"
#  - "StudyFuturePlan" (F2) renamed to "DroppedEducation"
# ---------------------------------------------------------------------------
$wsEducation = $wb.Worksheets.Item("Education")
$wsEducation.Range("B2").Value = "HighestQualification"
$wsEducation.Range("F2").Value = "DroppedEducation"

# ---------------------------------------------------------------------------
# Family sheet:
#  - "BrothersEmployed" (Q1) and "SistersEmployed" (R1) columns removed
# ---------------------------------------------------------------------------
$wsFamily = $wb.Worksheets.Item("Family")
$wsFamily.Range("Q1:R1").ClearContents()
